# Applies the "feat: new member as timestamp for box" edit:
#  1. Split the "27-02-25:" paragraph into a bold+underlined date run
#     ("27-02-25") followed by a plain ":" run.
#  2. Append a "(fatto)" run to the "Cercare di notificare senza comandi "
#     bullet paragraph.
#  3. Append six new paragraphs documenting the 06-03-25 entry at the end
#     of the document.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 1. "27-02-25:" -> bold+underlined "27-02-25" + plain ":"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(5)
$r1 = $p1.Range
$len1 = ($r1.Text).Length
$start1 = $r1.Start
$textRange1 = $d.Range($start1, $start1 + $len1 - 1)
$textRange1.Delete()

$xml1 = $pkgOpen + '<w:p><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="it-IT"/></w:rPr><w:t>27-02-25</w:t></w:r><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>:</w:t></w:r></w:p>' + $pkgClose
$d.Range($start1, $start1).InsertXML($xml1)

# ---------------------------------------------------------------------
# 2. "Cercare di notificare senza comandi " -> append "(fatto)" run
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(8)
$r2 = $p2.Range
$len2 = ($r2.Text).Length
$start2 = $r2.Start
$textRange2 = $d.Range($start2, $start2 + $len2 - 1)
$textRange2.Delete()

$xml2 = $pkgOpen + '<w:p><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">Cercare di notificare senza comandi </w:t></w:r><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>(fatto)</w:t></w:r></w:p>' + $pkgClose
$d.Range($start2, $start2).InsertXML($xml2)

# ---------------------------------------------------------------------
# 3. Append the new "06-03-25:" section at the end of the document
# ---------------------------------------------------------------------
$newParas = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="it-IT"/></w:rPr><w:t>06</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="it-IT"/></w:rPr><w:t>-0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="it-IT"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="it-IT"/></w:rPr><w:t>-25:</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">Al momento il BOT si avvia in collegamento con un chatID preciso </w:t></w:r><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">(messaggio di welcome) </w:t></w:r><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>e comunica l' + [char]0x2019 + 'allarme quando viene inviato da Thingsboard.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>Per semplicit' + [char]0x00E0 + ' diremo che abbiamo collegato due dispositivi (2 chat ID) associati a due cassette diverse.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>Implementato la possibilit' + [char]0x00E0 + ' di mandare l' + [char]0x2019 + 'ultimo accesso verso Thingsboard. Ora tramite comando lo si vorrebbe recuperare.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr></w:p>'

$xml3 = $pkgOpen + $newParas + $pkgClose
$endPos = $d.Content.End
$d.Range($endPos, $endPos).InsertXML($xml3)

Write-Output "edit applied"
